# Applies the "prompt_id 1" results blocks + highlights the existing
# "5-shot" rows, matching the commit "Changes to eval_results to evaluate
# with results with prompting / updates results.xlsx with prompt_id 1 run
# details."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Light-orange highlight used on every "5-shot" row (Accent2, Lighter 80%).
$highlightColor = 14083579   # 0xFBE5D6 packed as BGR for Range.Interior.Color

function Highlight-Row($rowNum) {
    $ws.Range("A$rowNum`:G$rowNum").Interior.Color = $highlightColor
}

# ---------------------------------------------------------------------
# 1) Highlight the pre-existing "5-shot" rows in all three tables.
# ---------------------------------------------------------------------
Highlight-Row 5
Highlight-Row 12
Highlight-Row 19

# ---------------------------------------------------------------------
# 2) Helper that writes one "Prompt_id 1" block: a bold title, the
#    standard 7-column header row, and a single highlighted 5-shot data
#    row whose Success/Failure/Unparseable counts feed the Pass rate /
#    Unparseable rate formulas (same shape as the "No prompting" tables
#    above, but scaled by 5 samples per problem).
# ---------------------------------------------------------------------
function Add-PromptBlock($titleRow, $success, $failure, $unparseable, $totalUnparseable) {
    $headerRow = $titleRow + 1
    $dataRow = $titleRow + 2

    $ws.Cells.Item($titleRow, 4).Value = "Prompt_id 1 - Results pass @ $k"
    $ws.Cells.Item($titleRow, 4).Font.Bold = $true
    $ws.Cells.Item($titleRow, 4).Locked = $false
    $ws.Cells.Item($titleRow, 5).Locked = $false
    $ws.Cells.Item($titleRow, 6).Locked = $false

    $ws.Cells.Item($headerRow, 1).Value = "Setting"
    $ws.Cells.Item($headerRow, 2).Value = "Success"
    $ws.Cells.Item($headerRow, 3).Value = "Failure"
    $ws.Cells.Item($headerRow, 4).Value = "Unparseable"
    $ws.Cells.Item($headerRow, 5).Value = "Pass rate"
    $ws.Cells.Item($headerRow, 6).Value = "Total Unparseable"
    $ws.Cells.Item($headerRow, 7).Value = "Unparseable rate"
    $ws.Range($ws.Cells.Item($headerRow, 1), $ws.Cells.Item($headerRow, 7)).Font.Bold = $true

    $ws.Cells.Item($dataRow, 1).Value = "5-shot"
    $ws.Cells.Item($dataRow, 2).Value = $success
    $ws.Cells.Item($dataRow, 3).Value = $failure
    $ws.Cells.Item($dataRow, 4).Value = $unparseable
    $ws.Cells.Item($dataRow, 5).Formula = "=B$dataRow/SUM(B$dataRow`:D$dataRow)"
    $ws.Cells.Item($dataRow, 5).NumberFormat = "0.00%"
    $ws.Cells.Item($dataRow, 6).Value = $totalUnparseable
    $ws.Cells.Item($dataRow, 7).Formula = "=F$dataRow/(SUM(B$dataRow`:D$dataRow)*5)"
    $ws.Cells.Item($dataRow, 7).NumberFormat = "0.00%"
    Highlight-Row $dataRow
}

# ---------------------------------------------------------------------
# 3) The three prompt_id-1 blocks (k=1, k=2, k=5), each separated by a
#    blank row, exactly as the author pasted them in below the original
#    no-prompting tables.
# ---------------------------------------------------------------------
$k = 1
Add-PromptBlock 22 13 35 55 55

$k = 2
Add-PromptBlock 26 22 46 35 115

$k = 5
Add-PromptBlock 30 43 53 7 274

# ---------------------------------------------------------------------
# 4) Leave the cursor where the author left it when they saved.
# ---------------------------------------------------------------------
$ws.Range("J21").Select()

Write-Host "done"
